$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B36").Value = "Projektilsystem"
$ws.Range("D36").Value = "4h"
$ws.Range("F36").Value = 43405
$ws.Range("F36").NumberFormat = $ws.Range("F33").NumberFormat

$ws.Range("B37").Value = "Projektilsystem machte extreme Performance Probleme, versucht zu beheben (kein Erfolg)"
$ws.Range("D37").Value = "4h"
$ws.Range("F37").Value = 43405
$ws.Range("F37").NumberFormat = $ws.Range("F33").NumberFormat

$ws.Range("B39").Value = "Umbau der ECS Architektur aus Optimierungsgründen"
$ws.Range("D39").Value = "5h"
$ws.Range("F39").Value = 43405
$ws.Range("F39").NumberFormat = $ws.Range("F33").NumberFormat

$ws.Range("D41").Select()
$excel.ActiveWindow.ScrollRow = 25
$excel.ActiveWindow.ScrollColumn = 1
